$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data table values for "Highest qualification level by age and gender" row
# Latest period (release date) and Next period (release date) columns
$ws.Range("C7").Value = "Jan 2024 - Dec 2024 (15/04/25)"
$ws.Range("D7").Value = "Jan 2025 - Dec 2025 (Apr 26)"

# Update the view state to match the saved workbook (scrolled down, D8 selected)
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D8").Select()

$wb.Save()
